# runtime update (2025-10-28 01:40:06)
# Refresh the "as_of_utc" timestamp (column AA) on every data row of the
# "Главные" and "Линейные" sheets, and update the stat columns for the
# three officials whose underlying numbers changed since the last pull.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "2025-10-27 16:27:31"
$newTimestamp = "2025-10-27 17:40:05"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}

# --- "Главные" sheet: Gashilov Viktor (row 9) ---
$wsMain = $wb.Worksheets.Item("Главные")
$wsMain.Cells.Item(9, 3).Value  = 19
$wsMain.Cells.Item(9, 4).Value  = 314
$wsMain.Cells.Item(9, 5).Value  = 169
$wsMain.Cells.Item(9, 6).Value  = 145
$wsMain.Cells.Item(9, 7).Value  = 16.53
$wsMain.Cells.Item(9, 8).Value  = 8.890000000000001
$wsMain.Cells.Item(9, 9).Value  = 7.63
$wsMain.Cells.Item(9, 10).Value = 82
$wsMain.Cells.Item(9, 11).Value = 70
$wsMain.Cells.Item(9, 21).Value = 1
$wsMain.Cells.Item(9, 23).Value = 16

# --- "Главные" sheet: Sidorenko Maksim (row 24) ---
$wsMain.Cells.Item(24, 3).Value  = 18
$wsMain.Cells.Item(24, 4).Value  = 294
$wsMain.Cells.Item(24, 5).Value  = 143
$wsMain.Cells.Item(24, 6).Value  = 151
$wsMain.Cells.Item(24, 7).Value  = 16.33
$wsMain.Cells.Item(24, 8).Value  = 7.94
$wsMain.Cells.Item(24, 9).Value  = 8.390000000000001
$wsMain.Cells.Item(24, 10).Value = 69
$wsMain.Cells.Item(24, 11).Value = 73
$wsMain.Cells.Item(24, 21).Value = 1
$wsMain.Cells.Item(24, 23).Value = 8

# --- "Линейные" sheet: Ivanov Yuriy (row 12) ---
$wsLines = $wb.Worksheets.Item("Линейные")
$wsLines.Cells.Item(12, 3).Value  = 17
$wsLines.Cells.Item(12, 4).Value  = 292
$wsLines.Cells.Item(12, 5).Value  = 150
$wsLines.Cells.Item(12, 6).Value  = 142
$wsLines.Cells.Item(12, 7).Value  = 17.18
$wsLines.Cells.Item(12, 8).Value  = 8.82
$wsLines.Cells.Item(12, 9).Value  = 8.35
$wsLines.Cells.Item(12, 10).Value = 75
$wsLines.Cells.Item(12, 11).Value = 66
$wsLines.Cells.Item(12, 21).Value = 1
$wsLines.Cells.Item(12, 23).Value = 4
